# Update the "想去人数" (want-to-go count) figures in column F
# for the "展览" (sheet1) and "全部类型" (sheet4) worksheets.
# These two sheets share the same underlying event rows (sheet4 simply
# has one extra row for an event that also lives on the "演出" sheet),
# so the row numbers differ by one from row 6 onward between them.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsAll  = $wb.Worksheets.Item("全部类型")

# 展览 sheet (rows 4-17)
$wsExpo.Range("F4").Value  = 82
$wsExpo.Range("F5").Value  = 9
$wsExpo.Range("F6").Value  = 548
$wsExpo.Range("F7").Value  = 7715
$wsExpo.Range("F8").Value  = 483
$wsExpo.Range("F9").Value  = 207
$wsExpo.Range("F11").Value = 685
$wsExpo.Range("F12").Value = 19
$wsExpo.Range("F13").Value = 28
$wsExpo.Range("F14").Value = 177
$wsExpo.Range("F15").Value = 32
$wsExpo.Range("F17").Value = 766

# 全部类型 sheet (rows 4-18, offset by +1 vs 展览 from row 6 onward)
$wsAll.Range("F4").Value  = 82
$wsAll.Range("F5").Value  = 9
$wsAll.Range("F7").Value  = 548
$wsAll.Range("F8").Value  = 7715
$wsAll.Range("F9").Value  = 483
$wsAll.Range("F10").Value = 207
$wsAll.Range("F12").Value = 685
$wsAll.Range("F13").Value = 19
$wsAll.Range("F14").Value = 28
$wsAll.Range("F15").Value = 177
$wsAll.Range("F16").Value = 32
$wsAll.Range("F18").Value = 766
